$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data for October 5, 2021 (row 8 / "5 tarikh")
$ws.Range("F8").Value = 1640

$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = 2
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 2

# Move selection/cursor as left by the editing session
[void]$ws.Range("L29").Select()
